# "code improvement and alignment"
# Adds a new "ConstantValues" worksheet holding a simple name/value table of
# workbook constants (currently just the delivery-charge limit), and makes
# it the active sheet/tab instead of "Categories".

$wb = $excel.ActiveWorkbook

# --- Add the new "ConstantValues" sheet, placed after the last existing sheet ---
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$constSheet = $wb.Worksheets.Add($null, $lastSheet)
$constSheet.Name = "ConstantValues"

# Column A is widened so the constant names are fully visible.
$constSheet.Columns.Item(1).ColumnWidth = 25.7

# Match the row-height convention used by the rest of the workbook.
$constSheet.Rows.Item(1).RowHeight = 12.8
$constSheet.Rows.Item(2).RowHeight = 12.8

# Header row.
$constSheet.Range("A1").Value = "Constant Name"
$constSheet.Range("B1").Value = "Value"
$constSheet.Range("A1:B1").Font.Bold = $true

# Data row - the one constant currently tracked in the sheet.
$constSheet.Range("A2").Value = "DELIVERYCHARGELIMIT"
$constSheet.Range("A2").Font.Name = "Times New Roman"
$constSheet.Range("B2").Value = 350

# The new sheet becomes the active tab (previously "Categories" was active).
$constSheet.Activate()
